$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "TestDataSheet"

# Resize the application window (best effort; engine re-derives bookViews on export)
$win = $excel.ActiveWindow
$win.Width = 14505
$win.Height = 6015

# Column widths for the newly introduced columns B, C, D
# (closest achievable values given the runtime's column-width quantization)
$ws.Columns.Item(2).ColumnWidth = 23
$ws.Columns.Item(3).ColumnWidth = 16.3
$ws.Columns.Item(4).ColumnWidth = 18

# Header row: shift TCID to A, Username/Password stay, SearchPassword -> SearchProduct
$ws.Range("A1").Value = "TCID"
$ws.Range("B1").Value = "Username"
$ws.Range("C1").Value = "Password"
$ws.Range("D1").Value = "SearchProduct"
$ws.Range("E1").Value = $null

# New data row
$ws.Range("A2").Value = "TC01"
$ws.Range("B2").Value = "tariq19ansari@gmail.com"
$ws.Range("C2").Value = "Qwerty123"
$ws.Range("D2").Value = "Julius Caesar"

# Hyperlink on the e-mail cell (adds the Hyperlink cell style automatically)
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:tariq19ansari@gmail.com")

# Final selection
$ws.Range("C10").Select() | Out-Null
